$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the 2007-2021 year-header row (row 4) and the percentage data row
# (row 5) with a new 2022 column (S), copying the adjacent cell's
# formatting (font/border/number-format/alignment) so the new cells match
# the look of the existing year columns.
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy($ws.Range("S5"))
$ws.Range("S5").Value = 76.1

# Match the saved selection state.
$ws.Range("P8").Select() | Out-Null
